# Update the "as_of_utc" timestamp column (AA) on the data sheets
# from "2025-11-07 02:49:21" to "2025-11-07 03:03:48".

$wb = $excel.ActiveWorkbook

$oldValue = "2025-11-07 02:49:21"
$newValue = "2025-11-07 03:03:48"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Range("AA$r")
        if ($cell.Value2 -eq $oldValue) {
            $cell.Value2 = $newValue
        }
    }
}
